# Realestate Update resale numbers 2025-02-04 22:22
# Appends a new data row (row 51) to the active sheet, mirroring the
# existing rows' layout: text columns A-D, numeric columns E-T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51

# Text columns (Date, Time, Weekday, Week) — use a leading apostrophe to
# force text interpretation (avoids Excel auto-converting "2025-02-04" to
# a date serial, "22:22:46" to a time serial, or "05" to the number 5),
# then reset the cell Style to "Normal" so no quotePrefix/number-format
# style gets attached to the cell (matches the plain, unstyled data rows).
function Set-TextCell($col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextCell 1 "2025-02-04"
Set-TextCell 2 "22:22:46"
Set-TextCell 3 "Tuesday"
Set-TextCell 4 "05"

# Numeric columns (city resale counts)
$ws.Cells.Item($row, 5).Value  = 125852
$ws.Cells.Item($row, 6).Value  = 141839
$ws.Cells.Item($row, 7).Value  = 166595
$ws.Cells.Item($row, 8).Value  = 157772
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 142082
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191043
$ws.Cells.Item($row, 14).Value = 115373
$ws.Cells.Item($row, 15).Value = 44763
$ws.Cells.Item($row, 16).Value = 28246
$ws.Cells.Item($row, 17).Value = 63067
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 39065
$ws.Cells.Item($row, 20).Value = -1
